$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the explanation text for the IPA row (D3): "IPM" -> "IPA" wording.
$ws.Range("D3").Value = "The Intergovernmental Personnel Act (IPA) Mobility Program allows agencies to receive temporary personnel assignments. This program is specifically focused on short-term engagements of non-Federal workers in the Federal space. Assignments may come to or from state and local governments, institutions of higher education, Indian tribal governments and other eligible organizations. IPA is useful to agencies looking for team members to bring innovation and perspective from outside the Federal government and provide a valuable experience for a non-Federal worker."

# The longer corrected text wraps to one additional line; match the
# recalculated row height for row 3.
$ws.Rows("3").RowHeight = 158.4

# Update the selected cell to D3, matching the saved view state.
$ws.Range("D3").Select()
